# #5: property boat&car done
# Bring the 汽車 (car) sheet in line with the other property sheets
# (土地/建物 etc.): turn row 1 into a real header row, add a "capacity"
# label for the existing engine-displacement column, and append the
# standard property metadata columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# Make room for the 7 new metadata columns (H:N) to the right of the
# existing data (A:G). Inserting one column at a time at "H" keeps
# pushing into empty space, so each new column simply inherits the
# formatting of the column immediately to its left (G) -- giving the
# header row (s=1) and the data rows (s=2) the same style as the rest
# of the table, exactly like Insert() does in real Excel.
$ws.Columns("H").Insert()
$ws.Columns("H").Insert()
$ws.Columns("H").Insert()
$ws.Columns("H").Insert()
$ws.Columns("H").Insert()
$ws.Columns("H").Insert()
$ws.Columns("H").Insert()

# Row 1 used to be a stray duplicate of row 2's data; replace it with the
# proper column headers (matching the 土地 / 建物 sheets).
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Row 2 (MAZDA3)
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2012-04-24"
$ws.Range("K2").Value = "李俊俋"
$ws.Range("L2").Value = 1738
$ws.Range("M2").Value = "tmp16861"
$ws.Range("N2").Value = 30

# Row 3 (MAZDAMPV)
$ws.Range("H3").Value = "land"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").Value = "2012-04-24"
$ws.Range("K3").Value = "李俊俋"
$ws.Range("L3").Value = 1738
$ws.Range("M3").Value = "tmp16861"
$ws.Range("N3").Value = 31

# Row 4 (三菱DELICA)
$ws.Range("H4").Value = "land"
$ws.Range("I4").Value = "normal"
$ws.Range("J4").Value = "2012-04-24"
$ws.Range("K4").Value = "李俊俋"
$ws.Range("L4").Value = 1738
$ws.Range("M4").Value = "tmp16861"
$ws.Range("N4").Value = 32
